# Auto-generated edit script: updates cryptos list figures (prices & 1h volume %)
# and fixes the BinanceUSD/Cosmos row order (rows 26-27 swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.258.02"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "1.651.22"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.92"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  +1.87%  "
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.256"
$ws.Range("E8").Value = "  +0.81%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.99"
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "1.881.91"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").Value = "1.642.57"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.67"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").Value = "27.233.30"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "220.66"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.84"
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("E22").Value = "  +6.61%  "
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.22"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.87"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.53"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("B27").Value = "BinanceUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.84"
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.40"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("D35").Value = "1.265.17"
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.547"
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.845"
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.45"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.20"
$ws.Range("E43").Value = "  +4.43%  "
$ws.Range("D44").Value = "1.791.92"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.17"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.74"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  +10.02%  "
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.68"
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0971"
$ws.Range("E51").Value = "  -0.72%  "
